# Add a new row to the "Recorded Classes Link" schedule table for
# "05th July" / "Data Sharing and List" / the session's YouTube link.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a new row at the end of the table; Word clones cell widths /
# run formatting (Times New Roman, 24 half-points) from the existing rows.
$newRow = $t.Rows.Add()

# --- Column 1: "05th July" (the "th" rendered as superscript) -----------
$cell1 = $newRow.Cells.Item(1)
$cell1Start = $cell1.Range.Start
$cell1.Range.Text = "05th July"
$thRange = $d.Range($cell1Start + 2, $cell1Start + 4)
$thRange.Font.Superscript = $true

# --- Column 2: topic ------------------------------------------------------
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Data Sharing and List"

# --- Column 3: YouTube hyperlink ------------------------------------------
$cell3 = $newRow.Cells.Item(3)
$d.Hyperlinks.Add($cell3.Range, "https://www.youtube.com/watch?v=J9i1zTUJnHw", `
    [Type]::Missing, [Type]::Missing, "https://www.youtube.com/watch?v=J9i1zTUJnHw") | Out-Null

Write-Output "Added row: 05th July / Data Sharing and List"
